# Update the COS470 schedule: shift all meeting dates back ~1 week, swap the
# Android/iOS week ordering ("android first"), and re-point the project-due
# notes (Project 1 <-> Android, Project 2 <-> iOS, Project 3 new URL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day (column B) shifts -------------------------------------------------
$ws.Range("B2").Value  = 19
$ws.Range("B3").Value  = 26
$ws.Range("B4").Value  = 2
$ws.Range("B5").Value  = 9
$ws.Range("B6").Value  = 16
$ws.Range("B7").Value  = 23
$ws.Range("B8").Value  = 2
$ws.Range("B9").Value  = 9
$ws.Range("B10").Value = 16
$ws.Range("B11").Value = 23
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 30
$ws.Range("B13").Value = 6
$ws.Range("B14").Value = 13
$ws.Range("B15").Value = 20
$ws.Range("B16").Value = 27
$ws.Range("B17").Value = 4
$ws.Range("B18").Value = 11

# --- Topic/reading content for weeks 3-8 (Android now comes before iOS) ----
$ws.Range("E3").Value = "Android (Kotlin) Basics Android Studio"
$ws.Range("F3").Value = "[Intro to Android](https://developer.android.com/guide/index.html), [Android Studio](https://developer.android.com/studio/intro/index.html), [Building Your First App](https://developer.android.com/training/basics/firstapp/index.html)"

$ws.Range("E4").Value = "Android (Kotlin) Activities and Intents"
$ws.Range("F4").Value = "[Ray Wenderlich Android Tutorials](https://www.raywenderlich.com/category/android), [Material Design To Do List Tutorial](http://dandroiddeveloper.github.io/list)"

$ws.Range("E5").Value = "Android (Kotlin)"
$ws.Range("F5").Value = "[Android Activities](https://developer.android.com/guide/components/activities/index.html), [React Native](http://facebook.github.io/react-native/), [React Native - Learn the Basics](http://facebook.github.io/react-native/docs/tutorial.html)"

$ws.Range("E6").Value = "iOS (Swift) Basics Xcode"
$ws.Range("F6").Value = "[Introduction to iOS 11, Xcode 9 and Swift 4](https://youtu.be/71pyOB4TPRE), [MVCs](https://youtu.be/l-2FaQTDYAw?list=PLPA-ayBrweUzGFmkT_W65z64MoGnKRZMq)"

$ws.Range("E7").Value = "iOS (Swift) Model-View-Controller"
$ws.Range("F7").Value = " [Start Developing iOS Apps](https://developer.apple.com/library/content/referencelibrary/GettingStarted/DevelopiOSAppsSwift/), [App Programming Guide for iOS](https://developer.apple.com/library/content/documentation/iPhone/Conceptual/iPhoneOSProgrammingGuide/Introduction/Introduction.html#//apple_ref/doc/uid/TP40007072-CH1-SW1), [Getting Started with iOS videos](https://training.apple.com/courses/Getting_Started_with_iOS_Development/training.html)"

$ws.Range("E8").Value = "iOS (Swift)"
$ws.Range("F8").Value = "[Swift Playgrounds](https://github.com/danielpi/Swift-Playgrounds)"

# --- Row heights follow the rotated content (rows 3-8) ----------------------
$ws.Rows.Item(3).RowHeight  = 119
$ws.Rows.Item(4).RowHeight  = 68
$ws.Rows.Item(5).RowHeight  = 119
$ws.Rows.Item(6).RowHeight  = 68
$ws.Rows.Item(7).RowHeight  = 187
$ws.Rows.Item(8).RowHeight  = 51

# --- Project-due notes (column G): Project 1 = Android, Project 2 = iOS ----
$ws.Range("G9").Value  = "[Project 2 - iOS Due](/projects/ios-tasklist)"
$ws.Range("G6").Value  = "[Project 1 - Android Due](/projects/android-tasklist)"
$ws.Range("G14").Value = "[Project 3 - Hybrid Due](/projects/hybrid-tasklist)"

# --- Selection / view ---------------------------------------------------
$ws.Range("F3").Select() | Out-Null
